# Rewrite the test-case rows on Sheet1 (the old "SSP adviser / personal client
# creation" scenario is replaced with the new "way2automation registration"
# scenario) and update the matching "continuation" rows on Sheet2, which hang
# off the same TestcaseID via the shared Steps/method_Name/testdata columns.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# Row 1 (header row: TestcaseType/TestcaseID/ScenarioName/Steps/method_Name/testdata)
# is unchanged, so only clear out the old body rows (2-8) before writing the
# new scenario's data - this also lets row 8 disappear entirely since the new
# scenario only needs rows 2-7.
$ws1.Range("A2:F8").ClearContents()

$ws1.Range("A2").Value = "Smoke"
$ws1.Range("B2").Value = "Client.001"
$ws1.Range("D2").Value = "User is on way2automationHomePage"
$ws1.Range("E2").Value = "Comman_Reusables.launchURL"

$ws1.Range("E3").Value = "Home_Page.validateUserIsOnHomePage"

$ws1.Range("D4").Value = "User clicks on the registration link"
$ws1.Range("E4").Value = "Home_Page.clickRegistrationLink"

$ws1.Range("D5").Value = "Registration pop up is displayed"
$ws1.Range("E5").Value = "Home_Page.validateRegistrationFormDisplayed"

$ws1.Range("D6").Value = "Register with a valid user"
$ws1.Range("E6").Value = "home_Page.registrationOnPopUp"
$ws1.Range("F6").Value = "Registration_testdata|Registration_Details|1"

$ws1.Range("E7").Value = "Registration_Page.registrationOnRegistrationPage"
$ws1.Range("F7").Value = "Registration_testdata|Registration_Details|1-2"

# ScenarioName is filled in last, which is why it lands at the tail of the
# shared-string table rather than right after TestcaseID.
$ws1.Range("C2").Value = "Client First Test Case"

# Sheet2 holds the trailing steps for the same Client.001 test case.
$ws2.Range("E13").Value = "Client_Page.createPersonalClient"
$ws2.Range("F13").Value = "Client_testdata|Client_Details|1|"
$ws2.Range("E14").Value = "Client_Page.clickSaveButtonOnSummaryPage"

# Restore Sheet2's own selection (unchanged) and move Sheet1's active cell to
# E19, matching the saved view state after the edit.
$ws2.Range("E14").Select()
$ws1.Range("E19").Select()
